$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8156.625
$ws.Range("J86").Value = 4501
$ws.Range("L86").Value = 4501
$ws.Range("N86").Value = -6747
$ws.Range("H88").Value = 7037.7
$ws.Range("J88").Value = 7580.778
$ws.Range("L88").Value = 7580.778
$ws.Range("N88").Value = -8392.778
$ws.Range("H89").Value = 8156.625
$ws.Range("J89").Value = 4501
$ws.Range("L89").Value = 22505
$ws.Range("N89").Value = -33737
$ws.Range("H91").Value = 7037.7
$ws.Range("J91").Value = 7580.778
$ws.Range("L91").Value = 7580.778
$ws.Range("N91").Value = -10388.778
$ws.Range("H106").Value = 2100
$ws.Range("I106").Value = 2100
$ws.Range("K106").Value = 2100
$ws.Range("M106").Value = -1469
$ws.Range("H125").Value = 28825
$ws.Range("J125").Value = 50650
$ws.Range("L125").Value = 455850
$ws.Range("N125").Value = -460770
$ws.Range("H137").Value = 2359
$ws.Range("I137").Value = 2037.7142
$ws.Range("K137").Value = 6113.142599999999
$ws.Range("M137").Value = -3563.142599999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 26318600
$ws.Range("I74").Value = 27029872
$ws.Range("K74").Value = 27029872
$ws.Range("M74").Value = -27028998
$ws.Range("H77").Value = 26318600
$ws.Range("I77").Value = 27029872
$ws.Range("K77").Value = 135149360
$ws.Range("M77").Value = -135144992
$ws.Range("H132").Value = 2502132.5
$ws.Range("I132").Value = 2633642
$ws.Range("J132").Value = 3450
$ws.Range("K132").Value = 7900926
$ws.Range("L132").Value = 10350
$ws.Range("M132").Value = -7898396
$ws.Range("N132").Value = -15410

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1460.3572
$ws.Range("I99").Value = 1420.0416
$ws.Range("J99").Value = 1702.25
$ws.Range("K99").Value = 1420.0416
$ws.Range("L99").Value = 1702.25
$ws.Range("M99").Value = 77.95839999999998
$ws.Range("N99").Value = -4698.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2744141.2
$ws.Range("J16").Value = 3500
$ws.Range("L16").Value = 3500
$ws.Range("N16").Value = -4074
$ws.Range("H31").Value = 13645.77
$ws.Range("J31").Value = 12699.583
$ws.Range("L31").Value = 12699.583
$ws.Range("N31").Value = -13289.583
$ws.Range("H34").Value = 13645.77
$ws.Range("J34").Value = 12699.583
$ws.Range("L34").Value = 12699.583
$ws.Range("N34").Value = -13103.583
$ws.Range("H58").Value = 38471124
$ws.Range("I58").Value = 45465236
$ws.Range("K58").Value = 45465236
$ws.Range("M58").Value = -45465033
$ws.Range("H105").Value = 2552875.2
$ws.Range("I105").Value = 4082922.5
$ws.Range("J105").Value = 2796.6667
$ws.Range("K105").Value = 4082922.5
$ws.Range("L105").Value = 2796.6667
$ws.Range("M105").Value = -4081175.5
$ws.Range("N105").Value = -6290.6667
$ws.Range("H113").Value = 2744141.2
$ws.Range("J113").Value = 3500
$ws.Range("L113").Value = 3500
$ws.Range("N113").Value = -7840
$ws.Range("H136").Value = 38471124
$ws.Range("I136").Value = 45465236
$ws.Range("K136").Value = 136395708
$ws.Range("M136").Value = -136393158

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 3085.5386
$ws.Range("I18").Value = 2711.3
$ws.Range("K18").Value = 8133.900000000001
$ws.Range("M18").Value = -7964.900000000001
$ws.Range("H32").Value = 143711.28
$ws.Range("I32").Value = 167496.67
$ws.Range("J32").Value = 999
$ws.Range("K32").Value = 502490.01
$ws.Range("L32").Value = 2997
$ws.Range("M32").Value = -502207.01
$ws.Range("N32").Value = -3563
$ws.Range("H61").Value = 149.33333
$ws.Range("I61").Value = 149.33333
$ws.Range("K61").Value = 447.99999
$ws.Range("M61").Value = -232.99999
$ws.Range("H75").Value = 302.125
$ws.Range("J75").Value = 308
$ws.Range("L75").Value = 924
$ws.Range("N75").Value = -2920
$ws.Range("H78").Value = 302.125
$ws.Range("J78").Value = 308
$ws.Range("L78").Value = 2772
$ws.Range("N78").Value = -12756
$ws.Range("H113").Value = 92217.73
$ws.Range("J113").Value = 1632.5
$ws.Range("L113").Value = 4897.5
$ws.Range("N113").Value = -9237.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1499.5
$ws.Range("I80").Value = 1999
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 1999
$ws.Range("L80").Value = 1000
$ws.Range("M80").Value = -1001
$ws.Range("N80").Value = -2996
$ws.Range("H83").Value = 1499.5
$ws.Range("I83").Value = 1999
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 9995
$ws.Range("L83").Value = 5000
$ws.Range("M83").Value = -5003
$ws.Range("N83").Value = -14984
$ws.Range("H111").Value = 34296.668
$ws.Range("J111").Value = 34296.668
$ws.Range("L111").Value = 34296.668
$ws.Range("N111").Value = -40430.668
$ws.Range("H113").Value = 74869.21000000001
$ws.Range("I113").Value = 94251.82000000001
$ws.Range("K113").Value = 94251.82000000001
$ws.Range("M113").Value = -92081.82000000001
$ws.Range("H122").Value = 129796.3
$ws.Range("I122").Value = 203662.17
$ws.Range("K122").Value = 610986.51
$ws.Range("M122").Value = -608536.51
$ws.Range("H126").Value = 2968
$ws.Range("J126").Value = 2971
$ws.Range("L126").Value = 8913
$ws.Range("N126").Value = -13853

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6143.6
$ws.Range("I61").Value = 6143.6
$ws.Range("K61").Value = 6143.6
$ws.Range("M61").Value = -5941.6
$ws.Range("H113").Value = 6143.6
$ws.Range("I113").Value = 6143.6
$ws.Range("K113").Value = 6143.6
$ws.Range("M113").Value = -3973.6
$ws.Range("H122").Value = 7224.8335
$ws.Range("I122").Value = 7899.75
$ws.Range("K122").Value = 23699.25
$ws.Range("M122").Value = -21249.25
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()
$ws.Range("H139").Value = 298994.75
$ws.Range("J139").Value = 298993
$ws.Range("L139").Value = 298993
$ws.Range("N139").Value = -309273

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1211.3636
$ws.Range("I96").Value = 912.25
$ws.Range("K96").Value = 912.25
$ws.Range("M96").Value = 460.75
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
